# feat: add 2022-Q1 data
#
# Before: sheets = [ "2021-Q4", "总计" ]
# After:  sheets = [ "2021-Q4", "2022-Q1", "总计" ]
#
# The existing "总计" (summary) sheet becomes the new "2022-Q1" sheet (keeps
# its sheetId / file / position-2 slot) and is repopulated with the fund
# holdings table for the new quarter. A duplicate of it is inserted right
# after (so it inherits the same sheetPr/header formatting), gets renamed
# back to "总计", and its summary table is updated to also list the freshly
# added "2022-Q1" row.

$wb = $excel.ActiveWorkbook

# Helper: write a value as genuine shared-string TEXT (never auto-coerced to
# a number), so strings like "0.10" keep their exact formatting instead of
# becoming 0.1. Uses a scratch cell far away, builds it via a text formula,
# then pastes the *value* (not the formula) into the destination.
function Set-TextValue {
    param($ws, [string]$cellRef, [string]$text)
    $scratch = $ws.Range("ZZ1000")
    $scratch.Formula = "=""" + $text + """"
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)   # xlPasteValues
    $scratch.Clear()
}

# ---------------------------------------------------------------------
# 1. Locate the current "总计" sheet; it sits right after "2021-Q4".
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$totalIndex = $wsTotal.Index

# ---------------------------------------------------------------------
# 2. Duplicate it, inserting the copy right after -> becomes position 3
#    with a fresh sheetId (matches the diff's new "总计", sheetId=3), and
#    already carries the original's sheetPr/outlinePr/header formatting.
# ---------------------------------------------------------------------
$wsTotal.Copy($null, $wsTotal)
$wsNewTotal = $wb.Worksheets.Item($totalIndex + 1)

# Extend the row-index column style down to the new row 3.
$wsNewTotal.Range("A2").Copy()
$wsNewTotal.Range("A3").PasteSpecial(-4122)      # xlPasteFormats

# Row 2: the newly added quarter goes first (overwrites the old "2021-Q4").
$wsNewTotal.Range("A2").Value = 0
Set-TextValue $wsNewTotal "B2" "2022-Q1"
$wsNewTotal.Range("C2").Value = 2
$wsNewTotal.Range("D2").Value = 0

# Row 3: the previous quarter, shifted down.
$wsNewTotal.Range("A3").Value = 1
Set-TextValue $wsNewTotal "B3" "2021-Q4"
$wsNewTotal.Range("C3").Value = 1
$wsNewTotal.Range("D3").Value = 0

# ---------------------------------------------------------------------
# 3. Turn the original "总计" sheet into the "2022-Q1" fund-holdings sheet.
# ---------------------------------------------------------------------

# Extend the bold/centered header style across the new columns E:H, and the
# row-index style down to the new row 3.
$wsTotal.Range("B1").Copy()
$wsTotal.Range("E1:H1").PasteSpecial(-4122)      # xlPasteFormats
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)         # xlPasteFormats

# New headers.
$wsTotal.Range("B1").Value = "基金代码"
$wsTotal.Range("C1").Value = "基金名称"
$wsTotal.Range("D1").Value = "基金规模"
$wsTotal.Range("E1").Value = "股票总仓位"
$wsTotal.Range("F1").Value = "仓位占比"
$wsTotal.Range("G1").Value = "持有市值(亿元)"
$wsTotal.Range("H1").Value = "仓位排名"

# Row 2.
$wsTotal.Range("A2").Value = 0
Set-TextValue $wsTotal "B2" "005126"
Set-TextValue $wsTotal "C2" "银河量化稳进混合"
Set-TextValue $wsTotal "D2" "0.10"
Set-TextValue $wsTotal "E2" "78.20"
Set-TextValue $wsTotal "F2" "1.93"
Set-TextValue $wsTotal "G2" "0.0019"
$wsTotal.Range("H2").Value = 7

# Row 3.
$wsTotal.Range("A3").Value = 1
Set-TextValue $wsTotal "B3" "001797"
Set-TextValue $wsTotal "C3" "华融新利灵活配置混合"
Set-TextValue $wsTotal "D3" "0.02"
Set-TextValue $wsTotal "E3" "48.66"
Set-TextValue $wsTotal "F3" "2.40"
Set-TextValue $wsTotal "G3" "0.0005"
$wsTotal.Range("H3").Value = 7

# ---------------------------------------------------------------------
# 4. Final renames -> "2021-Q4", "2022-Q1", "总计".
# ---------------------------------------------------------------------
$wsTotal.Name = "2022-Q1"
$wsNewTotal.Name = "总计"
